$wb = $excel.ActiveWorkbook

# The "想去人数" (interested count) figures were refreshed for a couple of
# events on both the "展览" sheet and the duplicated "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1679
    $ws.Range("F4").Value = 27
    $ws.Range("F6").Value = 461
}
